$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "factor" with "Nominal" and "numeric" with "Numeric" in column B (rows 2-78)
for ($r = 2; $r -le 78; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq "factor") {
        $cell.Value2 = "Nominal"
    } elseif ($cell.Value2 -eq "numeric") {
        $cell.Value2 = "Numeric"
    }
}

# Set column B width (explicit custom width matching column A's style of customWidth)
$ws.Columns.Item(2).ColumnWidth = 8

# Update the view: scroll so row 6 is at top, and select B12:B78 with active cell B12
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("B12:B78").Select()
